$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix the typo in the observable name: dash -> underscore before "phi_S"
$ws.Range("B2").Value = "A_LTcos(phi_h-phi_S)"

# Update the active cell selection on the sheet
$ws.Range("C4").Select()
